# Updated cryptos list on Mon Nov 11 06:39:49 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, avoiding Excel's auto-conversion of
# number-like strings (e.g. "618.50") into real numbers, and avoiding any
# lingering custom cell style/number-format after the write.
function Set-TextValue($range, $value) {
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# --- Rows 41-44: coins were re-ranked, so name/link/price/volume swapped
#     between rows (WhiteBITCoin <-> RenderToken, Stacks <-> dogwifhat)
Set-TextValue "B41" "RenderToken"
Set-TextValue "C41" "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue "D41" "5.91"
Set-TextValue "E41" "  +9.99%  "

Set-TextValue "B42" "WhiteBITCoin"
Set-TextValue "C42" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D42" "20.72"
Set-TextValue "E42" "  +3.81%  "

Set-TextValue "B43" "dogwifhat"
Set-TextValue "C43" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D43" "3.00"
Set-TextValue "E43" "  +20.28%  "

Set-TextValue "B44" "Stacks"
Set-TextValue "C44" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D44" "2.00"
Set-TextValue "E44" "  +12.23%  "

# --- Price (D) and Volume(1h) (E) refresh for the rest of the rows
Set-TextValue "D2" "81.200.88"
Set-TextValue "E2" "  +3.02%  "
Set-TextValue "D3" "3.146.74"
Set-TextValue "E3" "  -0.77%  "
Set-TextValue "E4" "  +0.19%  "
Set-TextValue "D5" "207.24"
Set-TextValue "E5" "  +0.92%  "
Set-TextValue "D6" "618.50"
Set-TextValue "E6" "  -1.79%  "
Set-TextValue "D7" "0.281"
Set-TextValue "E7" "  +24.51%  "
Set-TextValue "E8" "  +0.03%  "
Set-TextValue "D9" "0.577"
Set-TextValue "E9" "  -0.60%  "
Set-TextValue "D10" "3.146.45"
Set-TextValue "E10" "  -0.73%  "
Set-TextValue "D11" "0.575"
Set-TextValue "E11" "  -0.74%  "
Set-TextValue "D12" "0.0000253"
Set-TextValue "E12" "  +13.43%  "
Set-TextValue "E13" "  +0.28%  "
Set-TextValue "D14" "5.26"
Set-TextValue "E14" "  -2.80%  "
Set-TextValue "D15" "3.721.62"
Set-TextValue "E15" "  -0.82%  "
Set-TextValue "D16" "31.20"
Set-TextValue "E16" "  -0.94%  "
Set-TextValue "D17" "81.149.15"
Set-TextValue "E17" "  +2.98%  "
Set-TextValue "D18" "3.125.15"
Set-TextValue "E18" "  -1.41%  "
Set-TextValue "D19" "3.16"
Set-TextValue "E19" "  +11.51%  "
Set-TextValue "D20" "13.83"
Set-TextValue "E20" "  -3.78%  "
Set-TextValue "D21" "428.50"
Set-TextValue "E21" "  +0.34%  "
Set-TextValue "D22" "8.91"
Set-TextValue "E22" "  -4.80%  "
Set-TextValue "D23" "5.06"
Set-TextValue "E23" "  +3.00%  "
Set-TextValue "E24" "  +5.28%  "
Set-TextValue "D25" "5.15"
Set-TextValue "E25" "  +8.67%  "
Set-TextValue "D26" "3.304.41"
Set-TextValue "E26" "  -0.78%  "
Set-TextValue "D27" "75.31"
Set-TextValue "E27" "  -0.51%  "
Set-TextValue "D28" "10.77"
Set-TextValue "E28" "  -1.39%  "
Set-TextValue "E29" "  -0.02%  "
Set-TextValue "D30" "0.0000121"
Set-TextValue "E30" "  +5.88%  "
Set-TextValue "E31" "  +0.03%  "
Set-TextValue "D32" "8.90"
Set-TextValue "E32" "  +0.58%  "
Set-TextValue "D33" "559.57"
Set-TextValue "E33" "  +9.68%  "
Set-TextValue "D34" "0.148"
Set-TextValue "E34" "  +15.96%  "
Set-TextValue "D35" "1.47"
Set-TextValue "E35" "  +0.38%  "
Set-TextValue "D36" "0.150"
Set-TextValue "E36" "  +11.72%  "
Set-TextValue "E37" "  +0.29%  "
Set-TextValue "D38" "22.55"
Set-TextValue "E38" "  -1.35%  "
Set-TextValue "E39" "  +0.21%  "
Set-TextValue "D40" "0.404"
Set-TextValue "E40" "  +1.96%  "
Set-TextValue "D45" "159.55"
Set-TextValue "E45" "  -2.51%  "
Set-TextValue "D47" "186.30"
Set-TextValue "E47" "  -3.01%  "
Set-TextValue "D48" "45.28"
Set-TextValue "E48" "  +6.49%  "
Set-TextValue "D49" "1.31"
Set-TextValue "E49" "  +1.86%  "
Set-TextValue "E50" "  -5.48%  "
Set-TextValue "D51" "25.59"
Set-TextValue "E51" "  +4.21%  "
